$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values.
# A2, F2 and H2 look numeric/date-like ("915382", "03/03/2023", "500"), but in
# the source workbook every cell is stored as plain text (inlineStr), so we
# prefix those with a leading apostrophe to force Excel to keep them as text
# instead of auto-converting them into a number / date serial value.
$ws.Range("A2").Value = "'915382"
$ws.Range("C2").Value = "valaei"
$ws.Range("F2").Value = "'03/03/2023"
$ws.Range("G2").Value = "||animal"
$ws.Range("H2").Value = "'500"

# Row 3 (735554 / morteza / pashaei / ...) is removed entirely, so the
# dimension shrinks from A1:H3 to A1:H2.
$ws.Rows("3:3").Delete()
